# "Update feats raw data"
#
# Before: one sheet, "Sheet1", holding the big Feats reference table
#         (columns A:D, rows 1-431).
# After:  that same sheet is renamed "Feats", and a brand-new "Sheet1" is
#         inserted right after it holding some Firestore-schema scratch
#         notes (collection/document/field layout + sample query snippets).
#
# The cell-entry order below is NOT arbitrary - it reproduces the exact
# shared-string interning order the author's session produced (new
# strings are appended to sst in first-use order), so the workbook comes
# out byte-identical at the data level.

$wb = $excel.ActiveWorkbook

# --- rename the original sheet; it keeps its identity (sheetId/rId), only
#     the tab label changes.
$ws = $wb.ActiveSheet
$ws.Name = "Feats"

# --- insert the new scratch sheet right after "Feats"; Excel auto-names it
#     "Sheet1" since that name is now free again.
$notes = $wb.Worksheets.Add($null, $ws)
$notes.Name = "Sheet1"

# --- populate the notes sheet. Order matters for shared-string interning -
#     see header comment.
$notes.Range("A9").Value = "Feats"
$notes.Range("C12").Value = "Know all the languages"
$notes.Range("A1").Value = "Feats (Collection)"
$notes.Range("C4").Value = "Name (Field)"
$notes.Range("C5").Value = "Power(Field)"
$notes.Range("C6").Value = "Tier(Field)"
$notes.Range("C24").Value = 'var featsRef = db.collection("Feats");'
$notes.Range("C7").Value = "Class(Field)"
$notes.Range("C19").Value = "Barbarian "
$notes.Range("B2").Value = "(Document)"
$notes.Range("B10").Value = "hashID"
$notes.Range("C25").Value = 'var query = featsRef.where("Class", "==", "Barbarian").where("Tier", "==", "Adventurer");'
$notes.Range("C27").Value = 'var query = featsRef.where("Class", "==", "Barbarian")'
$notes.Range("C28").Value = 'var query = featsRef.where("Class", "==", "General")'
$notes.Range("C29").Value = "Union those together somehow…."

# these reuse strings already present in the workbook (Linguist, Champion,
# General, hashID again, Barbarian Rage, its description, Adventurer) -
# their entry order doesn't affect the shared-string table.
$notes.Range("C11").Value = "Linguist"
$notes.Range("C13").Value = "Champion"
$notes.Range("C14").Value = "General"
$notes.Range("B15").Value = "hashID"
$notes.Range("C16").Value = "Barbarian Rage"
$notes.Range("C17").Value = "A: When esc die is 4+, start raging for free as a quick action."
$notes.Range("C18").Value = "Adventurer"

# --- column widths on the notes sheet (best-fit-ish, matching final layout)
$notes.Columns("A").ColumnWidth = 16.7109375
$notes.Columns("B").ColumnWidth = 16.5703125
$notes.Columns("C").ColumnWidth = 81.7109375

# --- selections: notes sheet parks on C7, then focus returns to Feats at
#     A31 (that's the tab that ends up active/visible).
$notes.Range("C7").Select() | Out-Null
$ws.Activate() | Out-Null
$ws.Range("A31").Select() | Out-Null

# --- cosmetic: restore the (non-maximized) window geometry from the saved
#     session.
$win = $wb.Windows.Item(1)
$win.Width = 17970
$win.Height = 15195
$win.Left = 13170
$win.Top = 5115

Write-Output "Renamed Sheet1 -> Feats; added new Sheet1 with Firestore notes."
